$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: verdict FALSE -> INSUFFICIENT INFO, confidence 0.7 -> (blank, text-typed empty cell)
$ws.Range("B2").Value = "INSUFFICIENT INFO"
$ws.Range("C2").Value = "'"
$ws.Range("C2").ClearFormats()

# Row 3: verdict TRUE -> FALSE, confidence 0.9 -> 0
$ws.Range("B3").Value = "'FALSE"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = 0

# Row 4: verdict TRUE -> INSUFFICIENT INFO, confidence 0.9 -> (blank, text-typed empty cell)
$ws.Range("B4").Value = "INSUFFICIENT INFO"
$ws.Range("C4").Value = "'"
$ws.Range("C4").ClearFormats()

# Row 5: confidence 0.8 -> 1
$ws.Range("C5").Value = 1

# Row 8: verdict TRUE -> INSUFFICIENT INFO, confidence 0.8 -> (blank, text-typed empty cell)
$ws.Range("B8").Value = "INSUFFICIENT INFO"
$ws.Range("C8").Value = "'"
$ws.Range("C8").ClearFormats()

# Row 9: verdict TRUE -> INSUFFICIENT INFO, confidence 1 -> (blank, text-typed empty cell)
$ws.Range("B9").Value = "INSUFFICIENT INFO"
$ws.Range("C9").Value = "'"
$ws.Range("C9").ClearFormats()

# Row 11: confidence 0.9 -> 0.8
$ws.Range("C11").Value = 0.8

# Row 12: verdict INSUFFICIENT INFO -> TRUE, confidence (blank) -> 1
$ws.Range("B12").Value = "'TRUE"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = 1
